$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts all existing data,
# merged cells, column widths, and styles one column to the right.
$ws.Columns("A").Insert()

# New column A header cell (row 2): bold/filled/centered like the other
# merged section headers, then set its value.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "SyncChar"

# New column A body cells (rows 3-5): plain bordered/centered like the
# rest of the data rows, then set the row-3 value.
$ws.Range("B3").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$ws.Range("A3").Value = "M"

# Content updates elsewhere on the sheet (these land on I4 / G6 after the
# column insert shifted everything right by one column).
$ws.Range("I4").Value = "CAT Yellow"
$ws.Range("G6").Value = "H"

# Update the active selection to match the edited workbook (also drops
# the stale topLeftCell frozen-at-A2 setting).
[void]$ws.Range("A5").Select()
